$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the shared-string insertion order matching the target file:
# first-seen order must be ... Montpellier, Id du trajet, null ...

# Row 4: brand new row of data (uses "Montpellier" - introduces it first)
$ws.Range("A4").Value = 34
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 44114.416666666664
$ws.Range("E4").Value = 44114.416666666664
$ws.Range("F4").Value = 43545.635775462964
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = "Marseille"
$ws.Range("J4").Value = "Montpellier"
$ws.Range("K4").Value = 1

# Header row: rename "id du trajet :" to "Id du trajet"
$ws.Range("A1").Value = "Id du trajet"

# Row 3: add a string value in A3 ("null"), I3/J3 shift (Lyon/Bordeaux)
$ws.Range("A3").Value = "null"
$ws.Range("I3").Value = "Lyon"
$ws.Range("J3").Value = "Bordeaux"

# Row 2: add an id value in A2, keep rest, but I2/J2 shift (Paris/Marseille)
$ws.Range("A2").Value = 12
$ws.Range("I2").Value = "Paris"
$ws.Range("J2").Value = "Marseille"

# Match styling on new date cells to the existing date-format column style
# (columns D/E/F use a custom "yyyy/mm/dd hh:mm:ss" number format)
$ws.Range("D4").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"
$ws.Range("E4").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"
$ws.Range("F4").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"

# Column A width tweak
$ws.Columns.Item(1).ColumnWidth = 14

# Update active selection to match target
$ws.Range("D9").Select()
